$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 861
$ws.Range("B2").Value = 620
$ws.Range("C2").Value = 620
$ws.Range("D2").Value = 620
$ws.Range("E2").Value = 849
$ws.Range("F2").Value = 887
$ws.Range("G2").Value = 861
$ws.Range("H2").Value = 962
